$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.37427986513949918,
    -0.0099999992125638926,
    -0.0089999992027394171,
    -0.011999999808043782,
    -0.0059999992188455309,
    -0.0059999992033752392,
    -0.019999999063594842,
    -0.019999999062770613,
    -0.0059999992034907024,
    -0.0059999992063524132,
    -0.0044999992215153384,
    -0.0059999992080661535,
    -0.0059999992166126503,
    0.0643005261057219,
    -0.005999999217745966,
    -0.0059999992154824433,
    -0.0059999992124817325,
    -0.008999999182050189,
    -0.0089999992190601397,
    -0.072755274477626308,
    -0.0089999991796827494,
    -0.044576515953222096,
    -0.0089999991956366543,
    -0.041999998854906906,
    -0.041999998848523568,
    -0.0059999992023485049,
    -0.0059999991998429536,
    -0.0059999991895578475,
    -0.011999999123402105,
    -0.019999999041218075,
    0.03847974540193988,
    -0.020999999025989702,
    -0.031331791600192282
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
